$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column C ("solutions") - reorder the comma separated word lists
# ---------------------------------------------------------------------------
$ws.Range("C2").Value  = "answer, statement"
$ws.Range("C5").Value  = "plant, flower"
$ws.Range("C8").Value  = "geometry, pure_mathematics, mathematics"
$ws.Range("C12").Value = "celestial_body, star"
$ws.Range("C15").Value = "down, hill, ground"
$ws.Range("C28").Value = "sitting, furniture, seat, chair"
$ws.Range("C29").Value = "laboratory, lab"

# ---------------------------------------------------------------------------
# Column E ("relation") - updated relationship descriptions (related_to
# removed from the relation sets and remaining relations re-ordered)
# ---------------------------------------------------------------------------
$ws.Range("E2").Value  = 'question is "desires,  antonym, distinct_from" to answer | reply is "antonym,  synonym, " to answer | solution is " synonym" to answer'
$ws.Range("E4").Value  = 'antlers is "" to deer | doe is " " to deer | fawn is " " to deer'
$ws.Range("E5").Value  = 'bud is " " to flower | dandelion is " " to flower | petals is "" to flower'
$ws.Range("E6").Value  = 'colt is " " to horse | mare is "  part_of" to horse | unicorn is "" to horse'
$ws.Range("E12").Value = 'astronomy is " has_context" to star | moon is " distinct_from" to star | twinkle is "" to star'
$ws.Range("E13").Value = 'bait is "" to fish | pond is "used_for,  at_location" to fish | tuna is " " to fish'
$ws.Range("E15").Value = 'gravity is "" to down | low is "" to down | up is " antonym, distinct_from" to down'
$ws.Range("E17").Value = 'brawl is " " to fight | debate is "" to fight | soldier is "desires,  used_for" to fight'
$ws.Range("E19").Value = 'finger is "at_location,  part_of" to hand | glove is "" to hand | palm is " part_of" to hand'
$ws.Range("E21").Value = 'discuss is "" to talk | gossip is " " to talk | telephone is "" to talk'
$ws.Range("E28").Value = 'bench is "distinct_from, " to chair | sofa is "distinct_from,  " to chair | stool is "antonym,  synonym,  etymologically_" to chair'
$ws.Range("E33").Value = 'hand is " at_location, part_of" to finger | toe is "antonym, distinct_from, synonym, similar_to, " to finger | trigger is "" to finger'

# ---------------------------------------------------------------------------
# Column widths - widen columns E:H (5:8) to fit the refreshed content
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 25.83
$ws.Columns.Item(6).ColumnWidth = 24.33
$ws.Columns.Item(7).ColumnWidth = 34.5
$ws.Columns.Item(8).ColumnWidth = 25.5

# ---------------------------------------------------------------------------
# Reset the view: scroll back to the top-left corner and select A1 (clears
# the stale topLeftCell="B1" / selection at M28 left over from editing)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1").Select()
